# Auto-generated edit script applying the cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.254.95"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.495.69"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.58"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.45"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.35"
$ws.Range("E9").Value = "  +3.46%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("D12").Value = "4.090.90"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "3.496.57"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "64.345.50"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.61"
$ws.Range("E17").Value = "  -6.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.85"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.50"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "392.33"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.572"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "3.635.64"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.56"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.35"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  -7.19%  "
$ws.Range("D32").Value = "3.515.67"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").Value = "  -4.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.87"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "166.63"
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0779"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.99"
$ws.Range("E43").Value = "  -6.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.39"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "2.370.75"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.891"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0257"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("E51").Value = "  -1.51%  "
